# piv-opt.pptx: rename "Pivotal Query Optimizer" box to "GPORCA" (and shrink
# it to its new auto-fit size), merge the "Legacy "/"Optimizer" runs into a
# single "Legacy Optimizer" run, and refresh the cached "datetimeFigureOut"
# footer date across the slide master and every slide layout.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "Pivotal Query Optimizer" -> "GPORCA" -------------------------------
# This textbox lives inside "Group 70"; setting .Text re-runs PowerPoint's
# shrink-to-fit logic for the single remaining line, and then we nudge the
# box down to the position recorded for the shorter label.
$gporcaGroup = $s.Shapes.Item(6)
$gporcaBox = $gporcaGroup.GroupItems.Item(2)
$gporcaBox.TextFrame.TextRange.Text = "GPORCA"
$gporcaBox.Top = 199.9921259842520

# --- 2. Merge "Legacy " + "Optimizer" into one run --------------------------
$legacyGroup = $s.Shapes.Item(9)
$legacyBox = $legacyGroup.GroupItems.Item(1)
# Force a real content change first (to an unrelated placeholder string) so
# the two original runs get collapsed into a single new run, instead of
# being left untouched because the final text equals the old concatenation.
$legacyBox.TextFrame.TextRange.Text = "X"
$legacyBox.TextFrame.TextRange.Text = "Legacy Optimizer"

# --- 3. Refresh the footer date field everywhere ----------------------------
function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "6/6/2014") {
                # Same trick as above: change then set, so the cached field
                # text is actually rewritten rather than left as-is.
                $tr.Text = "X"
                $tr.Text = "2/2/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($i)
}
